$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.438.51'
$ws.Range('E2').Value = '  -2.20%  '

$ws.Range('D3').Value = '2.052.10'
$ws.Range('E3').Value = '  -0.56%  '

$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.12%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.663'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.70%  '

$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '53.88'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -9.24%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '57.98'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.64%  '

$ws.Range('E10').Value = '  -8.25%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0746'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.76%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.107'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.51%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.893'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.87%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.64'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -8.20%  '

$ws.Range('D15').Value = '2.350.36'
$ws.Range('E15').Value = '  -0.57%  '

$ws.Range('E16').Value = '  -9.11%  '

$ws.Range('D17').Value = '2.040.99'
$ws.Range('E17').Value = '  -0.83%  '

$ws.Range('D18').Value = '36.377.17'
$ws.Range('E18').Value = '  -2.27%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.53'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -12.33%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.01%  '

$ws.Range('D21').Value = '0.0₃0852'
$ws.Range('E21').Value = '  -6.44%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '236.83'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.18%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.23'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.35%  '

$ws.Range('E24').Value = '  +0.17%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.35'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.57%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.23'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.47%  '

$ws.Range('E27').Value = '  -5.81%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.52'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.37%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.98'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.74%  '

$ws.Range('E30').Value = '  -3.61%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.02'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -9.44%  '

$ws.Range('E32').Value = '  -0.51%  '

$ws.Range('E33').Value = '  -6.97%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0589'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.73%  '

$ws.Range('E35').Value = '  +0.12%  '

$ws.Range('E36').Value = '  +1.36%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0828'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.41%  '

$ws.Range('E38').Value = '  -8.26%  '

$ws.Range('E39').Value = '  -8.76%  '

$ws.Range('E40').Value = '  -7.74%  '

$ws.Range('E41').Value = '  -6.34%  '

$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.81'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -10.25%  '

$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.10'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.87%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '92.90'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -8.22%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0894'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -12.60%  '

$ws.Range('D46').Value = '1.372.78'
$ws.Range('E46').Value = '  +4.52%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.56'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -11.23%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.29'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.53%  '

$ws.Range('E49').Value = '  -1.72%  '

$ws.Range('E50').Value = '  -7.85%  '

$ws.Range('D51').Value = '2.238.10'
$ws.Range('E51').Value = '  -0.51%  '
